$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells whose styles never change, used as copy-style sources:
#   TextStyleSrc -> style 14 (General/text, right-aligned)
#   CountStyleSrc -> style 15 (#,##0, numeric count column)
#   PctStyleSrc   -> style 16 (#,##0.0, numeric percent column)
$textStyleSrc  = $ws.Range("A14")
$countStyleSrc = $ws.Range("F22")
$pctStyleSrc   = $ws.Range("H24")

function Copy-Style($srcCell, $dstCell) {
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)
}

function Set-Cell($ws, $addr, $newKind, $val, $catKind, $oldKind) {
    $cell = $ws.Range($addr)
    if ($newKind -eq "s") {
        if ($oldKind -ne "s") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $val
        if ($oldKind -ne "s") {
            Copy-Style $textStyleSrc $cell
        }
    } else {
        if ($oldKind -eq "s") {
            if ($catKind -eq "count") {
                Copy-Style $countStyleSrc $cell
            } else {
                Copy-Style $pctStyleSrc $cell
            }
        }
        $cell.Value2 = [double]$val
    }
}

Set-Cell $ws "D15" "s" "0" "count" "n"
Set-Cell $ws "E15" "s" "***.*" "pct" "n"
Set-Cell $ws "F15" "s" "0" "count" "n"
Set-Cell $ws "G15" "n" 3 "count" "n"
Set-Cell $ws "H15" "n" -100 "pct" "n"
Set-Cell $ws "C16" "n" 5 "count" "n"
Set-Cell $ws "D16" "n" 2 "count" "n"
Set-Cell $ws "E16" "n" 150 "pct" "n"
Set-Cell $ws "F16" "n" 17 "count" "n"
Set-Cell $ws "G16" "n" 23 "count" "n"
Set-Cell $ws "H16" "n" -26.086956521739 "pct" "n"
Set-Cell $ws "I16" "n" 119 "count" "n"
Set-Cell $ws "J16" "n" 82 "count" "n"
Set-Cell $ws "K16" "n" 45.121951219512 "pct" "n"
Set-Cell $ws "L16" "n" 120.37037037037 "pct" "n"
Set-Cell $ws "M16" "n" 75 "pct" "n"
Set-Cell $ws "N16" "n" -67.924528301886 "pct" "n"
Set-Cell $ws "C17" "n" 6 "count" "n"
Set-Cell $ws "D17" "n" 11 "count" "n"
Set-Cell $ws "E17" "n" -45.454545454545 "pct" "n"
Set-Cell $ws "F17" "n" 34 "count" "n"
Set-Cell $ws "G17" "n" 30 "count" "n"
Set-Cell $ws "H17" "n" 13.333333333333 "pct" "n"
Set-Cell $ws "I17" "n" 153 "count" "n"
Set-Cell $ws "J17" "n" 74 "count" "n"
Set-Cell $ws "K17" "n" 106.756756756757 "pct" "n"
Set-Cell $ws "L17" "n" 146.774193548387 "pct" "n"
Set-Cell $ws "M17" "n" 146.774193548387 "pct" "n"
Set-Cell $ws "N17" "n" 47.115384615384 "pct" "n"
Set-Cell $ws "C18" "n" 9 "count" "n"
Set-Cell $ws "D18" "n" 11 "count" "n"
Set-Cell $ws "E18" "n" -18.181818181818 "pct" "n"
Set-Cell $ws "G18" "n" 39 "count" "n"
Set-Cell $ws "H18" "n" 20.512820512820 "pct" "n"
Set-Cell $ws "I18" "n" 209 "count" "n"
Set-Cell $ws "J18" "n" 161 "count" "n"
Set-Cell $ws "K18" "n" 29.813664596273 "pct" "n"
Set-Cell $ws "L18" "n" 77.118644067796 "pct" "n"
Set-Cell $ws "M18" "n" 38.410596026490 "pct" "n"
Set-Cell $ws "N18" "n" -74.197530864197 "pct" "n"
Set-Cell $ws "C19" "n" 28 "count" "n"
Set-Cell $ws "D19" "n" 37 "count" "n"
Set-Cell $ws "E19" "n" -24.324324324324 "pct" "n"
Set-Cell $ws "F19" "n" 116 "count" "n"
Set-Cell $ws "G19" "n" 106 "count" "n"
Set-Cell $ws "H19" "n" 9.433962264150 "pct" "n"
Set-Cell $ws "I19" "n" 484 "count" "n"
Set-Cell $ws "J19" "n" 506 "count" "n"
Set-Cell $ws "K19" "n" -4.347826086956 "pct" "n"
Set-Cell $ws "L19" "n" 125.116279069767 "pct" "n"
Set-Cell $ws "M19" "n" 149.484536082474 "pct" "n"
Set-Cell $ws "N19" "n" 17.191283292978 "pct" "n"
Set-Cell $ws "C20" "n" 11 "count" "n"
Set-Cell $ws "D20" "n" 6 "count" "n"
Set-Cell $ws "E20" "n" 83.333333333333 "pct" "n"
Set-Cell $ws "F20" "n" 32 "count" "n"
Set-Cell $ws "G20" "n" 16 "count" "n"
Set-Cell $ws "I20" "n" 148 "count" "n"
Set-Cell $ws "J20" "n" 75 "count" "n"
Set-Cell $ws "K20" "n" 97.333333333333 "pct" "n"
Set-Cell $ws "L20" "n" 150.847457627119 "pct" "n"
Set-Cell $ws "M20" "n" 59.139784946236 "pct" "n"
Set-Cell $ws "N20" "n" -88.641596316193 "pct" "n"
Set-Cell $ws "C21" "n" 59 "count" "n"
Set-Cell $ws "D21" "n" 67 "count" "n"
Set-Cell $ws "E21" "n" -11.940298507462 "pct" "n"
Set-Cell $ws "F21" "n" 246 "count" "n"
Set-Cell $ws "G21" "n" 217 "count" "n"
Set-Cell $ws "H21" "n" 13.364055299539 "pct" "n"
Set-Cell $ws "I21" "n" 1125 "count" "n"
Set-Cell $ws "J21" "n" 910 "count" "n"
Set-Cell $ws "K21" "n" 23.626373626373 "pct" "n"
Set-Cell $ws "L21" "n" 117.601547388781 "pct" "n"
Set-Cell $ws "M21" "n" 96.335078534031 "pct" "n"
Set-Cell $ws "N21" "n" -62.636997675191 "pct" "n"
Set-Cell $ws "C22" "s" "0" "count" "n"
Set-Cell $ws "D22" "n" 1 "count" "s"
Set-Cell $ws "E22" "n" -100 "pct" "s"
Set-Cell $ws "G22" "n" 3 "count" "n"
Set-Cell $ws "H22" "n" -33.333333333333 "pct" "n"
Set-Cell $ws "J22" "n" 10 "count" "n"
Set-Cell $ws "K22" "n" -30 "pct" "n"
Set-Cell $ws "C23" "n" 1 "count" "n"
Set-Cell $ws "F23" "n" 3 "count" "n"
Set-Cell $ws "I23" "n" 6 "count" "n"
Set-Cell $ws "K23" "n" 50 "pct" "n"
Set-Cell $ws "L23" "n" 100 "pct" "n"
Set-Cell $ws "M23" "n" 100 "pct" "n"
Set-Cell $ws "C24" "n" 55 "count" "n"
Set-Cell $ws "D24" "n" 70 "count" "n"
Set-Cell $ws "E24" "n" -21.428571428571 "pct" "n"
Set-Cell $ws "F24" "n" 199 "count" "n"
Set-Cell $ws "G24" "n" 222 "count" "n"
Set-Cell $ws "H24" "n" -10.360360360360 "pct" "n"
Set-Cell $ws "I24" "n" 900 "count" "n"
Set-Cell $ws "J24" "n" 924 "count" "n"
Set-Cell $ws "K24" "n" -2.597402597402 "pct" "n"
Set-Cell $ws "L24" "n" 78.926441351888 "pct" "n"
Set-Cell $ws "M24" "n" 91.897654584221 "pct" "n"
Set-Cell $ws "C25" "n" 16 "count" "n"
Set-Cell $ws "D25" "n" 12 "count" "n"
Set-Cell $ws "E25" "n" 33.333333333333 "pct" "n"
Set-Cell $ws "F25" "n" 70 "count" "n"
Set-Cell $ws "G25" "n" 52 "count" "n"
Set-Cell $ws "H25" "n" 34.615384615384 "pct" "n"
Set-Cell $ws "I25" "n" 254 "count" "n"
Set-Cell $ws "J25" "n" 220 "count" "n"
Set-Cell $ws "K25" "n" 15.454545454545 "pct" "n"
Set-Cell $ws "L25" "n" 88.148148148148 "pct" "n"
Set-Cell $ws "M25" "n" 14.414414414414 "pct" "n"
Set-Cell $ws "D26" "s" "0" "count" "n"
Set-Cell $ws "E26" "s" "***.*" "pct" "n"
Set-Cell $ws "F26" "n" 4 "count" "n"
Set-Cell $ws "G26" "n" 7 "count" "n"
Set-Cell $ws "H26" "n" -42.857142857142 "pct" "n"
Set-Cell $ws "I26" "n" 21 "count" "n"
Set-Cell $ws "K26" "n" 31.25 "pct" "n"
Set-Cell $ws "L26" "n" 90.909090909090 "pct" "n"
Set-Cell $ws "C27" "s" "0" "count" "n"
Set-Cell $ws "D27" "n" 2 "count" "s"
Set-Cell $ws "E27" "n" -100 "pct" "s"
Set-Cell $ws "F27" "n" 4 "count" "n"
Set-Cell $ws "H27" "n" -33.333333333333 "pct" "n"
Set-Cell $ws "I27" "n" 36 "count" "n"
Set-Cell $ws "J27" "n" 31 "count" "n"
Set-Cell $ws "K27" "n" 16.129032258064 "pct" "n"
Set-Cell $ws "L27" "n" 89.473684210526 "pct" "n"
Set-Cell $ws "L28" "n" 200 "pct" "s"
Set-Cell $ws "L29" "n" 200 "pct" "s"

$excel.CutCopyMode = 0

# Update Volume Number text (A8): "16" -> "17"
$volCell = $ws.Range("A8")
$volText = $volCell.Text
$volCell.Characters($volText.IndexOf("16") + 1, 2).Text = "17"

# Update Report Covering the Week dates (C9)
$weekCell = $ws.Range("C9")
$weekText = $weekCell.Text
$weekCell.Characters($weekText.IndexOf("4/17/2023") + 1, 9).Text = "4/24/2023"
$weekText2 = $weekCell.Text
$weekCell.Characters($weekText2.IndexOf("4/23/2023") + 1, 9).Text = "4/30/2023"